$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# Row 2: H:1, Wholesale and retail trade (B unchanged)
$ws.Range("C2").Value = 130
$ws.Range("D2").Value = 369.8721498371335

# Row 3: H:1, Financial intermediation... -> Professional and technical services
$ws.Range("B3").Value = "Professional and technical services"
$ws.Range("C3").Value = 130
$ws.Range("D3").Value = 424.8439739413681

# Row 4: H:1, Community; social and personal services -> All other services
$ws.Range("B4").Value = "All other services"
$ws.Range("C4").Value = 130
$ws.Range("D4").Value = 383.4454397394137

# Row 5: H:2, Wholesale and retail trade (B unchanged)
$ws.Range("C5").Value = 260
$ws.Range("D5").Value = 369.8721498371335
$ws.Range("F5").Value = 2

# Row 6: H:2, Financial intermediation... -> Professional and technical services
$ws.Range("B6").Value = "Professional and technical services"
$ws.Range("C6").Value = 260
$ws.Range("D6").Value = 424.8439739413681
$ws.Range("F6").Value = 2

# Row 7: H:2, Community; social and personal services -> All other services
$ws.Range("B7").Value = "All other services"
$ws.Range("C7").Value = 260
$ws.Range("D7").Value = 383.4454397394137
$ws.Range("F7").Value = 2

# Row 8: H:3, Wholesale and retail trade (B unchanged)
$ws.Range("C8").Value = 450
$ws.Range("D8").Value = 369.8721498371335
$ws.Range("F8").Value = 3

# Row 9: H:3, Financial intermediation... -> Professional and technical services
$ws.Range("B9").Value = "Professional and technical services"
$ws.Range("C9").Value = 450
$ws.Range("D9").Value = 424.8439739413681
$ws.Range("F9").Value = 3

# Row 10: H:3, Community; social and personal services -> All other services
$ws.Range("B10").Value = "All other services"
$ws.Range("C10").Value = 450
$ws.Range("D10").Value = 383.4454397394137
$ws.Range("F10").Value = 3

# Row 11: HBET:3-6, Wholesale and retail trade (B unchanged)
$ws.Range("C11").Value = 900
$ws.Range("D11").Value = 369.8721498371335
$ws.Range("F11").Value = 5

# Row 12: HBET:3-6, Financial intermediation... -> Professional and technical services
$ws.Range("B12").Value = "Professional and technical services"
$ws.Range("C12").Value = 900
$ws.Range("D12").Value = 424.8439739413681
$ws.Range("F12").Value = 5

# Row 13: HBET:3-6, Community; social and personal services -> All other services
$ws.Range("B13").Value = "All other services"
$ws.Range("C13").Value = 900
$ws.Range("D13").Value = 383.4454397394137
$ws.Range("F13").Value = 5

# Row 14: HBET:4-7, Financial intermediation... -> Professional and technical services
$ws.Range("B14").Value = "Professional and technical services"
$ws.Range("C14").Value = 1200
$ws.Range("D14").Value = 556.5048859934853
$ws.Range("F14").Value = 5

# Row 15: HBET:4-7, Wholesale and retail trade -> Professional and technical services
$ws.Range("B15").Value = "Professional and technical services"
$ws.Range("C15").Value = 1200
$ws.Range("D15").Value = 369.8721498371335
$ws.Range("F15").Value = 5

# Row 16: HBET:4-7, Community; social and personal services -> All other services
$ws.Range("B16").Value = "All other services"
$ws.Range("C16").Value = 1200
$ws.Range("D16").Value = 383.4454397394137
$ws.Range("F16").Value = 5

# Row 17: HBET:8+, Financial intermediation... -> Professional and technical services
$ws.Range("B17").Value = "Professional and technical services"
$ws.Range("C17").Value = 3200
$ws.Range("D17").Value = 556.5048859934853
$ws.Range("F17").Value = 10

# Row 18: HBET:8+, Wholesale and retail trade -> Professional and technical services
$ws.Range("B18").Value = "Professional and technical services"
$ws.Range("C18").Value = 3200
$ws.Range("D18").Value = 369.8721498371335
$ws.Range("F18").Value = 10

# Row 19: HBET:8+, Community; social and personal services -> All other services
$ws.Range("B19").Value = "All other services"
$ws.Range("C19").Value = 3200
$ws.Range("D19").Value = 383.4454397394137
$ws.Range("F19").Value = 10
